# Weekly refresh of fruit/vegetable price data: reassign the
# Fecha/Volumen/Precio mínimo/Precio máximo/Precio promedio ponderado/
# Origen/Precio $/Kg values across the existing data rows (2-14).
# Row 7 is unaffected by this week's update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44232
$ws.Range("J2").Value = 250

# Row 3
$ws.Range("D3").Value = 44208

# Row 4
$ws.Range("D4").Value = 44204
$ws.Range("J4").Value = 430

# Row 5
$ws.Range("D5").Value = 44189
$ws.Range("J5").Value = 250

# Row 6
$ws.Range("D6").Value = 44230
$ws.Range("J6").Value = 250
$ws.Range("K6").Value = 5000
$ws.Range("M6").Value = 5500
$ws.Range("O6").Value = "Provincia de Quillota"
$ws.Range("P6").Value = 344

# Row 8
$ws.Range("D8").Value = 44186

# Row 9
$ws.Range("D9").Value = 44188
$ws.Range("J9").Value = 210
$ws.Range("L9").Value = 6000
$ws.Range("M9").Value = 5500
$ws.Range("O9").Value = "Provincia de Quillota"
$ws.Range("P9").Value = 344

# Row 10
$ws.Range("D10").Value = 44292
$ws.Range("J10").Value = 90
$ws.Range("K10").Value = 6000
$ws.Range("M10").Value = 6000
$ws.Range("O10").Value = "Región Metropolitana"
$ws.Range("P10").Value = 375

# Row 11
$ws.Range("D11").Value = 44187
$ws.Range("J11").Value = 160

# Row 12
$ws.Range("D12").Value = 44215

# Row 13
$ws.Range("D13").Value = 44251
$ws.Range("J13").Value = 120
$ws.Range("L13").Value = 5000
$ws.Range("M13").Value = 5000
$ws.Range("O13").Value = "Región Metropolitana"
$ws.Range("P13").Value = 312

# Row 14
$ws.Range("D14").Value = 44231
